# Aula 03 - Sincrono e Assincrono
# Edits anotacoes.docx per the target diff:
#  1. "stack" -> "Stack" (split run) and drop gramStart/gramEnd proofErr marks
#  2. "heap" -> "Heap" (split run) and drop gramStart/gramEnd proofErr marks
#  3. Move the _GoBack bookmark from the end of the ArrayList paragraph to the
#     end of the new "Estatico (static)" paragraph, and append a block of new
#     paragraphs (Estatico / Sincrono / Assincrono / Await / Task notes) right
#     after the blank paragraph that follows "List: permite colocar ...".

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    return -1
}

$stackIdx = Find-ParagraphIndex $d "stack:"
$heapIdx = Find-ParagraphIndex $d "heap:"
$listIdx = Find-ParagraphIndex $d "List:"

# --- 1 & 2: stack / heap paragraphs -------------------------------------
$stackXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:r><w:t>tack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: armazena valores menores, como INT, FLOAT, etc. Primeiro que entra é o último que sai.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item($stackIdx).Range.InsertXML($stackXml)

$heapXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>H</w:t></w:r><w:r><w:t>eap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: armazena valores maiores como objetos, referências, etc. Utiliza coletor de lixo.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item($heapIdx).Range.InsertXML($heapXml)

# --- 3: move the _GoBack bookmark, insert the new paragraphs -----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$blankIdx = $listIdx + 1
$insPoint = $d.Range($d.Paragraphs.Item($blankIdx).Range.End, $d.Paragraphs.Item($blankIdx).Range.End)
$newParasXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Estático (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>static</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>): não precisa instanciar, dar new()</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>Síncrono: executa linha por linha até o fim</w:t></w:r></w:p><w:p><w:r><w:t>Assíncrono (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>async</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>): tenta executar todas as operações possíveis</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Await</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: ‘bloqueia’ uma operação para o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>async</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> executar depois, por necessitar de resultados para funcionar melhor</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Await</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> retorna </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ou </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&lt;T&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> não tem valor de retorno</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&lt;T&gt; tem valor de retorno. Exemplo: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&gt;, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&lt;Produto&gt; </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Task.WhenAll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(t1, t2): espera as duas </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tasks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> terminarem para continuar</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insPoint.InsertXML($newParasXml)
